$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update M159:M164 (period_flux) from 10 to 24
# ---------------------------------------------------------------------------
foreach ($r in 159..164) {
    $ws.Cells.Item($r, 13).Value = 24
}

# ---------------------------------------------------------------------------
# Append the fabricio2019_thermal extraction (rows 165-167)
# ---------------------------------------------------------------------------

# Row 165
$ws.Cells.Item(165, 1).Value  = 1
$ws.Cells.Item(165, 2).Value  = 1
$ws.Cells.Item(165, 3).Value  = "fabricio2019_thermal"
$ws.Cells.Item(165, 4).Value  = "Figure 2 "
$ws.Cells.Item(165, 5).Value  = 0
$ws.Cells.Item(165, 6).Value  = 1
$ws.Cells.Item(165, 7).Value  = 1
$ws.Cells.Item(165, 8).Value  = 25
$ws.Cells.Item(165, 9).Font.Color = 0
$ws.Cells.Item(165, 9).NumberFormat = "0.0"
$ws.Cells.Item(165, 9).Value  = 20
$ws.Cells.Item(165, 10).Value = 15
$ws.Cells.Item(165, 11).Value = 25
$ws.Cells.Item(165, 12).Formula = "=K165-J165"
$ws.Cells.Item(165, 13).Value = 24
$ws.Cells.Item(165, 14).Value = "RMR"
$ws.Cells.Item(165, 15).Value = "ml O2 * kg^-1 *h^-1"
$ws.Cells.Item(165, 16).Value = 13.153153153153101
$ws.Cells.Item(165, 17).Value = 11.891891891891801
$ws.Cells.Item(165, 18).Value = 8.2749416080000007
$ws.Cells.Item(165, 19).Value = 3.0697364030697787
$ws.Cells.Item(165, 20).Value = 13
$ws.Cells.Item(165, 21).Value = 13
$ws.Cells.Item(165, 22).Value = 1
$ws.Cells.Item(165, 23).Value = 1
$ws.Cells.Item(165, 24).Value = "Crotalus"
$ws.Cells.Item(165, 25).Value = "durissus"
$ws.Cells.Item(165, 26).Value = 1
$ws.Cells.Item(165, 27).Value = 2
$ws.Cells.Item(165, 28).Value = 2

# Row 166
$ws.Cells.Item(166, 1).Value  = 1
$ws.Cells.Item(166, 2).Value  = 2
$ws.Cells.Item(166, 3).Value  = "fabricio2019_thermal"
$ws.Cells.Item(166, 4).Value  = "Figure 2 "
$ws.Cells.Item(166, 5).Value  = 0
$ws.Cells.Item(166, 6).Value  = 1
$ws.Cells.Item(166, 7).Value  = 1
$ws.Cells.Item(166, 8).Value  = 25
$ws.Cells.Item(166, 9).Font.Color = 0
$ws.Cells.Item(166, 9).NumberFormat = "0.0"
$ws.Cells.Item(166, 9).Value  = 25
$ws.Cells.Item(166, 10).Value = 20
$ws.Cells.Item(166, 11).Value = 30
$ws.Cells.Item(166, 12).Formula = "=K166-J166"
$ws.Cells.Item(166, 13).Value = 24
$ws.Cells.Item(166, 14).Value = "RMR"
$ws.Cells.Item(166, 15).Value = "ml O2 * kg^-1 *h^-1"
$ws.Cells.Item(166, 16).Value = 17.218181818181801
$ws.Cells.Item(166, 17).Value = 17.799999999999901
$ws.Cells.Item(166, 18).Value = 8.9427609427609625
$ws.Cells.Item(166, 19).Value = 3.0168350168349622
$ws.Cells.Item(166, 20).Value = 13
$ws.Cells.Item(166, 21).Value = 13
$ws.Cells.Item(166, 22).Value = 1
$ws.Cells.Item(166, 23).Value = 1
$ws.Cells.Item(166, 24).Value = "Crotalus"
$ws.Cells.Item(166, 25).Value = "durissus"
$ws.Cells.Item(166, 26).Value = 1
$ws.Cells.Item(166, 27).Value = 2
$ws.Cells.Item(166, 28).Value = 2

# Row 167
$ws.Cells.Item(167, 1).Value  = 1
$ws.Cells.Item(167, 2).Value  = 3
$ws.Cells.Item(167, 3).Value  = "fabricio2019_thermal"
$ws.Cells.Item(167, 4).Value  = "Figure 2 "
$ws.Cells.Item(167, 5).Value  = 0
$ws.Cells.Item(167, 6).Value  = 1
$ws.Cells.Item(167, 7).Value  = 1
$ws.Cells.Item(167, 8).Value  = 25
$ws.Cells.Item(167, 9).Font.Color = 0
$ws.Cells.Item(167, 9).NumberFormat = "0.0"
$ws.Cells.Item(167, 9).Value  = 30
$ws.Cells.Item(167, 10).Value = 25
$ws.Cells.Item(167, 11).Value = 35
$ws.Cells.Item(167, 12).Formula = "=K167-J167"
$ws.Cells.Item(167, 13).Value = 24
$ws.Cells.Item(167, 14).Value = "RMR"
$ws.Cells.Item(167, 15).Value = "ml O2 * kg^-1 *h^-1"
$ws.Cells.Item(167, 16).Value = 35.698529411764603
$ws.Cells.Item(167, 17).Value = 28.75
$ws.Cells.Item(167, 18).Value = 13.53485838779963
$ws.Cells.Item(167, 19).Value = 6.8627450980392553
$ws.Cells.Item(167, 20).Value = 13
$ws.Cells.Item(167, 21).Value = 13
$ws.Cells.Item(167, 22).Value = 1
$ws.Cells.Item(167, 23).Value = 1
$ws.Cells.Item(167, 24).Value = "Crotalus"
$ws.Cells.Item(167, 25).Value = "durissus"
$ws.Cells.Item(167, 26).Value = 1
$ws.Cells.Item(167, 27).Value = 2
$ws.Cells.Item(167, 28).Value = 2

# ---------------------------------------------------------------------------
# Update the view: select the last-entered cell (mirrors the saved selection)
# ---------------------------------------------------------------------------
$null = $ws.Range("S167").Select()

Write-Host "edit complete"
